$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 352 (pushes existing rows 352-436 down to 353-437)
$ws.Rows.Item(352).Insert()

# Populate the new row 352 with the new market record.
# Columns that are constant across the whole dataset for this sheet
# (A, B, C, E, F, G, H, I, N, Q, R) are carried over from the row that
# used to occupy 352 (now at 353), i.e. they keep the same values.
$ws.Cells.Item(352, 1).Value = 10
$ws.Cells.Item(352, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(352, 3).Value = "La Araucanía"
$ws.Cells.Item(352, 4).Value = 44782
$ws.Cells.Item(352, 5).Value = 9
$ws.Cells.Item(352, 6).Value = 100112008
$ws.Cells.Item(352, 7).Value = "Coliflor"
$ws.Cells.Item(352, 8).Value = "Sin especificar"
$ws.Cells.Item(352, 9).Value = "Primera"
$ws.Cells.Item(352, 10).Value = 200
$ws.Cells.Item(352, 11).Value = 1500
$ws.Cells.Item(352, 12).Value = 1500
$ws.Cells.Item(352, 13).Value = 1500
$ws.Cells.Item(352, 14).Value = '$/unidad'
$ws.Cells.Item(352, 15).Value = "Región Metropolitana"
$ws.Cells.Item(352, 16).Value = 1500
$ws.Cells.Item(352, 17).Value = 1
$ws.Cells.Item(352, 18).Value = "Hortaliza"
